$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Plan1 is the sheet shown as active/selected in the diff
$ws.Activate()

# Update specialty names (B2:B5) with the DATASUS dataset values
$ws.Range("B2").Value = "Clínico"
$ws.Range("B3").Value = "Cardilogia"
$ws.Range("B4").Value = "Pediatria"
$ws.Range("B5").Value = "Infectologia"

# Widen column B to fit the new, longer text (target stored width 22.44140625 chars)
$ws.Columns.Item(2).ColumnWidth = 21.65

# Update the selected cell shown in the sheet view
$ws.Range("J8").Select()

$wb.Save()
